$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 4635
    $ws.Range("F9").Value = 914
    $ws.Range("F11").Value = 1075
    $ws.Range("F13").Value = 586
    $ws.Range("F15").Value = 15
    $ws.Range("F16").Value = 268
}
